# Edit: add a new weekly entry (Primera/Segunda) for Acelga at the top of the
# historical data block (rows 128-129), pushing all subsequent rows down by 2
# (old row 128 -> 130, ..., old row 209 -> 211). This matches the commit
# "Fruta / hortaliza, semanal" which adds the latest week's prices while
# keeping history intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the first data row of this block (row 128).
# Everything currently at/after row 128 shifts down by 2 rows automatically,
# which reproduces rows 130-211 of the target sheet.
$ws.Range("A128:A129").EntireRow.Insert()

# New week's date (2021-08-30 == serial 44438)
$newDate = Get-Date -Year 2021 -Month 8 -Day 30 -Hour 0 -Minute 0 -Second 0

# Row 128: "Primera" quality entry for the new week
$ws.Cells.Item(128, 1).Value  = 8
$ws.Cells.Item(128, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(128, 3).Value  = "Coquimbo"
$ws.Cells.Item(128, 4).Value  = $newDate
$ws.Cells.Item(128, 5).Value  = 4
$ws.Cells.Item(128, 6).Value  = 100112009
$ws.Cells.Item(128, 7).Value  = "Acelga"
$ws.Cells.Item(128, 8).Value  = "Sin especificar"
$ws.Cells.Item(128, 9).Value  = "Primera"
$ws.Cells.Item(128, 10).Value = 3200
$ws.Cells.Item(128, 11).Value = 450
$ws.Cells.Item(128, 12).Value = 500
$ws.Cells.Item(128, 13).Value = 475
$ws.Cells.Item(128, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(128, 15).Value = "Provincia del Elqu$([char]0x00ED)"
$ws.Cells.Item(128, 16).Value = 238
$ws.Cells.Item(128, 17).Value = 2
$ws.Cells.Item(128, 18).Value = "Hortaliza"

# Row 129: "Segunda" quality entry for the new week
$ws.Cells.Item(129, 1).Value  = 8
$ws.Cells.Item(129, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(129, 3).Value  = "Coquimbo"
$ws.Cells.Item(129, 4).Value  = $newDate
$ws.Cells.Item(129, 5).Value  = 4
$ws.Cells.Item(129, 6).Value  = 100112009
$ws.Cells.Item(129, 7).Value  = "Acelga"
$ws.Cells.Item(129, 8).Value  = "Sin especificar"
$ws.Cells.Item(129, 9).Value  = "Segunda"
$ws.Cells.Item(129, 10).Value = 1540
$ws.Cells.Item(129, 11).Value = 350
$ws.Cells.Item(129, 12).Value = 400
$ws.Cells.Item(129, 13).Value = 375
$ws.Cells.Item(129, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(129, 15).Value = "Provincia del Elqu$([char]0x00ED)"
$ws.Cells.Item(129, 16).Value = 188
$ws.Cells.Item(129, 17).Value = 2
$ws.Cells.Item(129, 18).Value = "Hortaliza"

# Make sure the date cells keep the same date/time number format used
# throughout column D.
$ws.Cells.Item(128, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(129, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
